# CIERRE 8 ENE 2022
# Weekly payroll sheet ("RECIBOS NOMINA 2020") closing update: new amounts
# for this week's entries plus a refresh of the TODAY()-based closing date.
# Every SUM()/reference formula in the sheet recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (rows 3-7) ----------------------------------------------------
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 1833
$ws.Range("K4").Value = 433

# --- Block 2 (rows 20-26) ---------------------------------------------------
$ws.Range("K21").Value = 2380

# --- Block 3 (rows 36-41) ---------------------------------------------------
$ws.Range("K39").Value = 1250

# Refresh the "today" date chain: C14=TODAY() -> I14=C14 -> C32=I14 ->
# I32=C32 -> C48=C32 -> I48=C48 -> C65=I48. Re-asserting the root formula
# forces every downstream reference cell to recalc to the new closing day.
$ws.Range("C14").Formula = "=TODAY()"

# Leave the view where this week's closing entries were made.
$null = $ws.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("I62").Select()
